# "paises.xlsx" update: refresh the timestamp, re-rank a few countries
# whose case counts changed order, and write the new case figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header timestamp (row 1)
$ws.Range("A1").Value = "Datos actualizados a 16 de Abril de 2020 a las 01:22"

# Countries whose ranking moved (new counts changed the sort order) -
# update the country name shown in column A for each affected row.
$ws.Range("A107").Value = "San Marino"
$ws.Range("A108").Value = "Reunion"
$ws.Range("A109").Value = "Estado de Palestina"
$ws.Range("A148").Value = "Islas Caimanes"
$ws.Range("A149").Value = "Liberia"
$ws.Range("A150").Value = "Cabo Verde"
$ws.Range("A151").Value = "Polinesia Francesa"
$ws.Range("A152").Value = "Uganda"
$ws.Range("A153").Value = "Guyana"
$ws.Range("A154").Value = "Bahamas"
$ws.Range("A155").Value = "San Martin (Parte Holandesa)"
$ws.Range("A156").Value = "Guinea Ecuatorial"

# Updated figures: Casos totales / Nuevos casos / Casos activos /
# Recuperados / Casos criticos / Muertes hoy / Muertes
$ws.Range("B4").Value = 643296
$ws.Range("C4").Value = 29410
$ws.Range("D4").Value = 48567
$ws.Range("E4").Value = 566239
$ws.Range("G4").Value = 2443
$ws.Range("H4").Value = 28490

$ws.Range("B14").Value = 28610
$ws.Range("C14").Value = 3348
$ws.Range("E14").Value = 12827
$ws.Range("G14").Value = 225
$ws.Range("H14").Value = 1757

$ws.Range("B20").Value = 14350
$ws.Range("C20").Value = 124
$ws.Range("E20").Value = 5859

$ws.Range("B36").Value = 6301
$ws.Range("C36").Value = 190
$ws.Range("E36").Value = 5316

$ws.Range("B91").Value = 654
$ws.Range("C91").Value = 16
$ws.Range("D91").Value = 146
$ws.Range("E91").Value = 502

$ws.Range("B92").Value = 641
$ws.Range("C92").Value = 5
$ws.Range("D92").Value = 83
$ws.Range("E92").Value = 550

$ws.Range("B97").Value = 493
$ws.Range("C97").Value = 1
$ws.Range("D97").Value = 272
$ws.Range("E97").Value = 212
$ws.Range("F97").Value = 13
$ws.Range("G97").Value = 1
$ws.Range("H97").Value = 9

$ws.Range("B107").Value = 393
$ws.Range("C107").Value = 21
$ws.Range("D107").Value = 53
$ws.Range("E107").Value = 304
$ws.Range("F107").Value = 15
$ws.Range("H107").Value = 36

$ws.Range("B108").Value = 391
$ws.Range("C108").Value = 0
$ws.Range("D108").Value = 237
$ws.Range("E108").Value = 154
$ws.Range("F108").Value = 3
$ws.Range("H108").Value = 0

$ws.Range("B109").Value = 374
$ws.Range("C109").Value = 66
$ws.Range("D109").Value = 63
$ws.Range("E109").Value = 309
$ws.Range("F109").Value = 0
$ws.Range("H109").Value = 2

$ws.Range("B148").Value = 60
$ws.Range("C148").Value = 6
$ws.Range("D148").Value = 6
$ws.Range("E148").Value = 53
$ws.Range("F148").Value = 3
$ws.Range("H148").Value = 1

$ws.Range("B149").Value = 59
$ws.Range("C149").Value = 0
$ws.Range("D149").Value = 4
$ws.Range("E149").Value = 49
$ws.Range("H149").Value = 6

$ws.Range("B150").Value = 56
$ws.Range("C150").Value = 45
$ws.Range("D150").Value = 1
$ws.Range("E150").Value = 54
$ws.Range("F150").Value = 0
$ws.Range("H150").Value = 1

$ws.Range("D151").Value = 0
$ws.Range("E151").Value = 55
$ws.Range("F151").Value = 1

$ws.Range("C152").Value = 0
$ws.Range("D152").Value = 12
$ws.Range("E152").Value = 43
$ws.Range("F152").Value = 0
$ws.Range("H152").Value = 0

$ws.Range("B153").Value = 55
$ws.Range("C153").Value = 8
$ws.Range("D153").Value = 8
$ws.Range("E153").Value = 41
$ws.Range("F153").Value = 5
$ws.Range("H153").Value = 6

$ws.Range("C154").Value = 4
$ws.Range("D154").Value = 6
$ws.Range("F154").Value = 1
$ws.Range("H154").Value = 8

$ws.Range("B155").Value = 53
$ws.Range("C155").Value = 1
$ws.Range("D155").Value = 5
$ws.Range("E155").Value = 39
$ws.Range("F155").Value = 2
$ws.Range("H155").Value = 9

$ws.Range("B156").Value = 51
$ws.Range("C156").Value = 10
$ws.Range("D156").Value = 4
$ws.Range("E156").Value = 47
$ws.Range("F156").Value = 0
$ws.Range("H156").Value = 0
